$wb = $excel.ActiveWorkbook
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# Typography sheet: the "Default" typography's Wildcard Characters needs a "."
# added so the decimal point in the D3 ADC voltage reading (e.g. "3.301") can
# be rendered by the font.
$wsTypography.Range("G4").Value = "."

# Translation sheet: SingleUseId1 format string gets a trailing "v" unit
# suffix, and SingleUseId2 now holds the D3 Pin's ADC voltage reading instead
# of the old placeholder value.
$wsTranslation.Range("F4").Value = "ADC value = <value>v"

# F5 must stay text ("3.301"), not get auto-converted to a number. Temporarily
# force a Text format, assign the value, then restore the default style so no
# numeric coercion occurs and no extra formatting lingers on the cell.
$cellF5 = $wsTranslation.Range("F5")
$cellF5.NumberFormat = "@"
$cellF5.Value = "3.301"
$cellF5.Style = "Normal"
